# Update van de tijdsbestedingen
# - B4 becomes a formula 5+2+5 (was a literal 3)
# - B5 becomes a formula =5 (was empty)
# - E2 (SUM) and E3 (AVERAGE) recalc automatically as a consequence
# - Active selection moves from B5 to B6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Formula = "=5+2+5"
$ws.Range("B5").Formula = "=5"

$ws.Range("B6").Select()
